$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13: August 2019
$ws.Range("A13").Value = "August 2019"
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat

$ws.Range("B13").Value = "https://myemail.constantcontact.com/News-From-The-Forest---August.html?soid=1102494320279&aid=ePsj_Z-h1SI"
$ws.Hyperlinks.Add($ws.Range("B13"), "https://myemail.constantcontact.com/News-From-The-Forest---August.html?soid=1102494320279&aid=ePsj_Z-h1SI")
$ws.Range("B13").Style = $ws.Range("B12").Style

# Row 14: September 2019
$ws.Range("A14").Value = "September 2019"
$ws.Range("A14").NumberFormat = $ws.Range("A12").NumberFormat

$ws.Range("B14").Value = "https://myemail.constantcontact.com/News-From-The-Forest---September.html?soid=1102494320279&aid=TJ5FtWA1WX4"
$ws.Hyperlinks.Add($ws.Range("B14"), "https://myemail.constantcontact.com/News-From-The-Forest---September.html?soid=1102494320279&aid=TJ5FtWA1WX4")
$ws.Range("B14").Style = $ws.Range("B12").Style

$ws.Range("B14").Select()
